$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the table with two more rows (2025-09-02) for both stations, matching the
# layout/number-formats of the two rows directly above (A: date style, F: integer style).
$ws.Range("A64:F65").Copy()
$ws.Range("A66:F67").PasteSpecial(-4122)

# Row 66: 四方坪站
$ws.Range("A66").Value = 45902
$ws.Range("B66").Value = "四方坪站"
$ws.Range("C66").Value = 9665.7099999999991
$ws.Range("D66").Value = 7793.54
$ws.Range("E66").Value = 3320.37
$ws.Range("F66").Value = 421

# Row 67: 高岭站
$ws.Range("A67").Value = 45902
$ws.Range("B67").Value = "高岭站"
$ws.Range("C67").Value = 4496.63
$ws.Range("D67").Value = 3637.26
$ws.Range("E67").Value = 1148.75
$ws.Range("F67").Value = 159

[void]$ws.Range("G69").Select()
